$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update employee data (FirstName, LastName, Age, EmployeeCategory)
$ws.Range("A2").Value = "Javlon"
$ws.Range("B2").Value = "Fayziyev"
$ws.Range("C2").Value = 22
$ws.Range("D2").Value = 3

$ws.Range("A3").Value = "Qosimjon"
$ws.Range("B3").Value = "Berdiyev"
$ws.Range("C3").Value = 19
$ws.Range("D3").Value = 11

$ws.Range("A4").Value = "Farangiz"
$ws.Range("C4").Value = 15
$ws.Range("D4").Value = 13

$ws.Range("A5").Value = "Muhammad Rizo"
$ws.Range("B5").Value = "Keldiyev"
$ws.Range("C5").Value = 18
$ws.Range("D5").Value = 5

$ws.Range("B4").Value = "Gadoyeva"

# Update selection to E7
$ws.Range("E7").Select()
